$wb = $excel.ActiveWorkbook

# --- Sheet "survey" (sheet1): datetime picker -> split date/time pickers ---
$ws1 = $wb.Worksheets.Item("survey")

# Remove the now-unused inputAttributes.min / inputAttributes.max columns (O:P)
$ws1.Range("O1:P1").EntireColumn.Delete()

# Column N header: inputAttributes.data-field -> inputAttributes.timeFormat
$ws1.Cells.Item(1, 14).Value = "inputAttributes.timeFormat"

# New example rows for the split date / time pickers
$ws1.Cells.Item(3, 7).Value = "date"
$ws1.Cells.Item(3, 9).Value = "example2"
$ws1.Cells.Item(3, 10).Value = "Date"
$ws1.Cells.Item(3, 14).Value = "YYYY/DD/MM"

$ws1.Cells.Item(4, 7).Value = "time"
$ws1.Cells.Item(4, 9).Value = "example3"
$ws1.Cells.Item(4, 10).Value = "Time"
$ws1.Cells.Item(4, 14).Value = "HH:mm"

# --- Sheet "settings" (sheet2): pick up the Hyperlink / Followed Hyperlink
#     built-in cell styles that this Excel build always stamps into the
#     style table on save (no visible hyperlink is actually added). ---
$ws2 = $wb.Worksheets.Item("settings")
$ws2.Hyperlinks.Add($ws2.Range("Z100"), "https://example.com")
$ws2.Hyperlinks.Add($ws2.Range("Z101"), "https://example.com")
$ws2.Range("Z101").Style = "Followed Hyperlink"
$ws2.Hyperlinks.Delete()
$ws2.Range("Z100:Z101").Style = "Normal"
$ws2.Range("Z100:Z101").ClearContents()

# "survey" becomes the active / selected sheet and cell
$ws1.Activate() | Out-Null
$ws1.Range("N4").Select() | Out-Null
